# Financials update: insert a new "latest period" column (D) into the LTM
# sheet's three statement blocks (Income Statement, Balance Sheet, Cash Flow
# Statement), shifting the existing D:K data right to E:L, and populate the
# new column with the newest period's figures. Also correct two data points
# that were restated (Retained Earnings FY2012 and Capital Expenditures
# history) as part of the same update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column at D - this shifts all existing D:K content
#    (and their cell styles) right to E:L, and extends the used range.
$ws.Columns("D:D").Insert()

# 2) The freshly inserted column D cells default to the generic style, not
#    the date / number style used by the rest of the table. Clone the
#    formatting from column E (which holds the old "D" formatting after the
#    shift) onto the new column D so D matches its neighbours exactly.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Populate column D with the new period's values, row by row.
#    (Row numbers below correspond to the three "Period Ending" blocks:
#    Income Statement rows 7-35, Balance Sheet rows 38-77, Cash Flow rows
#    80-102.)
$newColumnD = @{
    7  = 43465
    8  = 9895500
    9  = 7962800
    10 = 1932600
    12 = "NA"
    13 = 0
    14 = 0
    15 = 0
    17 = 9136800
    18 = 758600
    20 = -115600
    21 = 643000
    22 = 346000
    23 = 297000
    24 = 83800
    25 = 0
    26 = 213200
    27 = 181900
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 115600
    33 = 181900
    34 = 0
    35 = 181900
    38 = 43465
    41 = 622400
    42 = 800500
    43 = 1234600
    44 = 279300
    45 = 369500
    46 = 3306400
    47 = 64800
    48 = 9953400
    49 = 3735100
    50 = 0
    51 = 0
    52 = 507100
    53 = 0
    54 = 17566800
    57 = 1048000
    58 = 1397200
    59 = 3123600
    60 = 5568800
    61 = 5864600
    62 = 2386700
    63 = 0
    64 = 0
    65 = 0
    66 = 13899900
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = 3259300
    73 = 0
    74 = 0
    75 = 0
    76 = 3666800
    77 = 0
    80 = 43465
    81 = 181900
    83 = 0
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = 1516900
    91 = -660700
    92 = 0
    93 = 0
    94 = -358400
    96 = -72600
    97 = 0
    98 = 0
    99 = 0
    100 = -1052200
    101 = -166700
    102 = -60400
}

foreach ($row in ($newColumnD.Keys | Sort-Object)) {
    $ws.Cells.Item($row, 4).Value = $newColumnD[$row]
}

# 4) Two historical data points were restated alongside the new column:
#    - Retained Earnings (row 72): the FY2012 figure (now column J after the
#      shift) changes from 3,748,400 to 6,283,500.
#    - Capital Expenditures (row 91): the FY2017-FY2013 figures (now columns
#      E:I after the shift) are restated.
$ws.Cells.Item(72, 10).Value = 6283500

$ws.Cells.Item(91, 5).Value = -403700
$ws.Cells.Item(91, 6).Value = -694400
$ws.Cells.Item(91, 7).Value = -1569700
$ws.Cells.Item(91, 8).Value = -1440400
$ws.Cells.Item(91, 9).Value = -1381800
